$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values would otherwise be auto-detected as
# numbers by Excel, so they remain text (matching the original string type).
$textCells = @("D4", "D5", "D6", "D7", "D10", "D13", "D17", "D19", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D33", "D35", "D36", "D38", "D40", "D41", "D42", "D44", "D45")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "51.471.86"
$ws.Range("E2").Value = "  -1.11%  "

# Row 3
$ws.Range("D3").Value = "2.930.09"
$ws.Range("E3").Value = "  -2.53%  "

# Row 4
$ws.Range("D4").Value = "0.999"

# Row 5
$ws.Range("D5").Value = "373.89"
$ws.Range("E5").Value = "  +5.54%  "

# Row 6
$ws.Range("D6").Value = "103.10"
$ws.Range("E6").Value = "  -3.56%  "

# Row 7
$ws.Range("D7").Value = "0.542"
$ws.Range("E7").Value = "  -2.97%  "

# Row 8
$ws.Range("E8").Value = "  -0.26%  "

# Row 9
$ws.Range("E9").Value = "  -4.08%  "

# Row 10
$ws.Range("D10").Value = "36.94"
$ws.Range("E10").Value = "  -2.96%  "

# Row 11
$ws.Range("E11").Value = "  -0.69%  "

# Row 12
$ws.Range("E12").Value = "  -2.40%  "

# Row 13
$ws.Range("D13").Value = "18.34"
$ws.Range("E13").Value = "  -3.51%  "

# Row 14
$ws.Range("D14").Value = "3.391.27"
$ws.Range("E14").Value = "  -2.49%  "

# Row 15
$ws.Range("E15").Value = "  -3.39%  "

# Row 16
$ws.Range("D16").Value = "2.935.48"
$ws.Range("E16").Value = "  -1.83%  "

# Row 17
$ws.Range("D17").Value = "0.927"
$ws.Range("E17").Value = "  -8.47%  "

# Row 18
$ws.Range("D18").Value = "51.428.42"
$ws.Range("E18").Value = "  -1.34%  "

# Row 19
$ws.Range("D19").Value = "3.41"
$ws.Range("E19").Value = "  +0.42%  "

# Row 20
$ws.Range("E20").Value = "  -1.84%  "

# Row 21
$ws.Range("D21").Value = "12.92"
$ws.Range("E21").Value = "  -4.66%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0946"
$ws.Range("E22").Value = "  -2.69%  "

# Row 23
$ws.Range("D23").Value = "68.27"
$ws.Range("E23").Value = "  -1.21%  "

# Row 24
$ws.Range("D24").Value = "261.98"
$ws.Range("E24").Value = "  -0.64%  "

# Row 25
$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  +0.98%  "

# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "4.14"
$ws.Range("E26").Value = "  -4.84%  "

# Row 27
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.168"
$ws.Range("E27").Value = "  -5.57%  "

# Row 28
$ws.Range("E28").Value = "  -0.03%  "

# Row 29
$ws.Range("D29").Value = "25.76"
$ws.Range("E29").Value = "  -4.59%  "

# Row 30
$ws.Range("D30").Value = "7.30"
$ws.Range("E30").Value = "  -2.23%  "

# Row 31
$ws.Range("D31").Value = "6.94"
$ws.Range("E31").Value = "  +6.85%  "

# Row 32
$ws.Range("E32").Value = "  -5.26%  "

# Row 33
$ws.Range("D33").Value = "9.81"
$ws.Range("E33").Value = "  -3.72%  "

# Row 34
$ws.Range("E34").Value = "  -3.25%  "

# Row 35
$ws.Range("D35").Value = "51.07"
$ws.Range("E35").Value = "  -0.03%  "

# Row 36
$ws.Range("D36").Value = "34.00"
$ws.Range("E36").Value = "  -5.62%  "

# Row 37
$ws.Range("E37").Value = "  +0.43%  "

# Row 38
$ws.Range("D38").Value = "0.0424"
$ws.Range("E38").Value = "  -2.76%  "

# Row 39
$ws.Range("E39").Value = "  -9.36%  "

# Row 40
$ws.Range("D40").Value = "16.98"
$ws.Range("E40").Value = "  -3.59%  "

# Row 41
$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  -9.09%  "

# Row 42
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  -7.12%  "

# Row 43
$ws.Range("E43").Value = "  -2.39%  "

# Row 44
$ws.Range("D44").Value = "123.53"
$ws.Range("E44").Value = "  -0.54%  "

# Row 45
$ws.Range("D45").Value = "21.66"
$ws.Range("E45").Value = "  -5.95%  "

# Row 46
$ws.Range("E46").Value = "  -5.23%  "

# Row 47
$ws.Range("E47").Value = "  +10.90%  "

# Row 48
$ws.Range("D48").Value = "2.021.01"
$ws.Range("E48").Value = "  -4.72%  "

# Row 49
$ws.Range("E49").Value = "  -1.74%  "

# Row 50
$ws.Range("E50").Value = "  -5.23%  "

# Row 51
$ws.Range("D51").Value = "3.211.71"
$ws.Range("E51").Value = "  -2.79%  "
